# The first paragraph currently contains four runs:
#   "Thank you for your response. Now, please answer the following question:"
#   <br/>
#   <br/>
#   "1. How do you think technology has changed the way we communicate with each other?"
# It needs to collapse into a single run with the new prompt text, while
# keeping the paragraph mark (and its formatting) intact.

$d = $word.ActiveDocument
$para = $d.Paragraphs.First
$r = $para.Range
# Exclude the paragraph-mark character at the end of the range so we
# only replace the run content, not the paragraph break itself.
$r.End = $r.End - 1
$r.Text = "It seems like your response is unclear. Please provide your complete answer for the IELTS Speaking task on ""Describe a memorable meal you had."" Once you do, I will assess your response based on the criteria given."
